# MacManes_summary.docx edit
#
# "desert-animals amazing ability" -> "desert-animal's amazing ability"
# (the straight apostrophe is auto-corrected to a typographic right
# single quote, exactly as real Word's AutoFormat/AutoCorrect would do)
# and the document's "_GoBack" last-edit-position bookmark is relocated
# to sit right after the newly inserted apostrophe (i.e. immediately
# before the "s" that now starts "s amazing ability ...") - matching
# where Word leaves _GoBack after an in-place edit at that spot.

$d = $word.ActiveDocument

# 1) Fix "desert-animals" -> "desert-animal's" (curly apostrophe via
#    AutoCorrect) using Find/Replace scoped to the whole story.
$found = $d.Content.Find.Execute(
    "desert-animals amazing", $true, $false, $false, $false, $false,
    $true, 1, $false, "desert-animal's amazing", 2)

if (-not $found) {
    throw "Could not find target phrase 'desert-animals amazing' to replace."
}

# 2) Re-seat the _GoBack bookmark immediately before "s amazing ability"
#    (i.e. right after the apostrophe we just inserted). Re-adding a
#    bookmark with an existing name moves it rather than duplicating it.
$full = $d.Content.Text
$pos = $full.IndexOf("s amazing ability to avert")

if ($pos -lt 0) {
    throw "Could not locate insertion point for the _GoBack bookmark."
}

$target = $d.Range($pos, $pos)
$d.Bookmarks.Add("_GoBack", $target)
